$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed ticker cells (rows 2-18) to match the refreshed screener output
$ws.Range("B2").Value = 'NSE:AJMERA'
$ws.Range("C2").Value = 'NSE:AARTISURF'
$ws.Range("D2").Value = 'NSE:CHOLAFIN'
$ws.Range("E2").Value = 'NSE:AUROPHARMA'
$ws.Range("F2").Value = 'NSE:AXISBANK'
$ws.Range("B3").Value = 'NSE:APTUS'
$ws.Range("C3").Value = 'NSE:ADFFOODS'
$ws.Range("D3").Value = 'NSE:EICHERMOT'
$ws.Range("F3").Value = 'NSE:BRITANNIA'
$ws.Range("B4").Value = 'NSE:ASAHIINDIA'
$ws.Range("C4").Value = 'NSE:BANG'
$ws.Range("D4").Value = 'NSE:GODREJCP'
$ws.Range("F4").Value = 'NSE:EICHERMOT'
$ws.Range("B5").Value = 'NSE:ASTRAZEN'
$ws.Range("C5").Value = 'NSE:BEDMUTHA'
$ws.Range("D5").Value = 'NSE:NAVINFLUOR'
$ws.Range("F5").Value = 'NSE:GODREJPROP'
$ws.Range("C6").Value = 'NSE:BHAGERIA'
$ws.Range("D6").Value = 'NSE:NMDC'
$ws.Range("F6").Value = 'NSE:HINDUNILVR'
$ws.Range("B7").Value = 'NSE:AXISBNKETF'
$ws.Range("C7").Value = 'NSE:BLS'
$ws.Range("D7").Value = 'NSE:PIDILITIND'
$ws.Range("F7").Value = 'NSE:ICICIBANK'
$ws.Range("B8").Value = 'NSE:BRITANNIA'
$ws.Range("C8").Value = 'NSE:GREAVESCOT'
$ws.Range("D8").Value = 'NSE:PIIND'
$ws.Range("F8").Value = 'NSE:JUBLFOOD'
$ws.Range("B9").Value = 'NSE:CONCORDBIO'
$ws.Range("C9").Value = 'NSE:GULFPETRO'
$ws.Range("F9").Value = 'NSE:MARICO'
$ws.Range("B10").Value = 'NSE:EICHERMOT'
$ws.Range("C10").Value = 'NSE:HCL-INSYS'
$ws.Range("F10").Value = 'NSE:NESTLEIND'
$ws.Range("B11").Value = 'NSE:FIVESTAR'
$ws.Range("C11").Value = 'NSE:KANORICHEM'
$ws.Range("B12").Value = 'NSE:GODREJPROP'
$ws.Range("C12").Value = 'NSE:MOLDTKPAC'
$ws.Range("B13").Value = 'NSE:HDFCPVTBAN'
$ws.Range("C13").Value = 'NSE:MONTECARLO'
$ws.Range("B14").Value = 'NSE:HINDUNILVR'
$ws.Range("C14").Value = 'NSE:PANACEABIO'
$ws.Range("B15").Value = 'NSE:IVZINGOLD'
$ws.Range("C15").Value = 'NSE:RATNAMANI'
$ws.Range("B16").Value = 'NSE:JUBLFOOD'
$ws.Range("B17").Value = 'NSE:LODHA'
$ws.Range("B18").Value = 'NSE:MARICO'

# Clear cells that no longer have a ticker in the refreshed data
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()

# Remove the now-unused trailing rows (19-27); this also shrinks the sheet dimension to A1:F18
$ws.Rows("19:27").Delete()

